# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2841.2856
$ws.Range("J17").Value = 2841.2856
$ws.Range("L17").Value = 8523.856800000001
$ws.Range("N17").Value = -8859.856800000001
$ws.Range("H40").Value = 3093.5557
$ws.Range("J40").Value = 2930.2
$ws.Range("L40").Value = 2930.2
$ws.Range("N40").Value = -3280.2
$ws.Range("H64").Value = 5876.32
$ws.Range("I64").Value = 4367.25
$ws.Range("J64").Value = 7269.3076
$ws.Range("K64").Value = 4367.25
$ws.Range("L64").Value = 7269.3076
$ws.Range("M64").Value = -4119.25
$ws.Range("N64").Value = -7765.3076
$ws.Range("H67").Value = 5876.32
$ws.Range("I67").Value = 4367.25
$ws.Range("J67").Value = 7269.3076
$ws.Range("K67").Value = 4367.25
$ws.Range("L67").Value = 7269.3076
$ws.Range("M67").Value = -3509.25
$ws.Range("N67").Value = -8985.3076
$ws.Range("H132").Value = 2411.4443
$ws.Range("I132").Value = 1500.5
$ws.Range("J132").Value = 9699
$ws.Range("K132").Value = 4501.5
$ws.Range("L132").Value = 29097
$ws.Range("M132").Value = -1971.5
$ws.Range("N132").Value = -34157
$ws.Range("H137").Value = 11247.78
$ws.Range("I137").Value = 4630.56
$ws.Range("K137").Value = 13891.68
$ws.Range("M137").Value = -11341.68
$ws.Range("H138").Value = 3105.0308
$ws.Range("I138").Value = 2497.5715
$ws.Range("J138").Value = 4965.375
$ws.Range("K138").Value = 7492.7145
$ws.Range("L138").Value = 14896.125
$ws.Range("M138").Value = -2352.7145
$ws.Range("N138").Value = -25176.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2000.4375
$ws.Range("I45").Value = 2003.5
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2003.5
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1626.5
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 5411.857
$ws.Range("I61").Value = 4590.974
$ws.Range("K61").Value = 4590.974
$ws.Range("M61").Value = -4378.974
$ws.Range("H74").Value = 2638.4
$ws.Range("I74").Value = 1629.5
$ws.Range("J74").Value = 4345.769
$ws.Range("K74").Value = 1629.5
$ws.Range("L74").Value = 4345.769
$ws.Range("M74").Value = -755.5
$ws.Range("N74").Value = -6093.769
$ws.Range("H77").Value = 2638.4
$ws.Range("I77").Value = 1629.5
$ws.Range("J77").Value = 4345.769
$ws.Range("K77").Value = 8147.5
$ws.Range("L77").Value = 21728.845
$ws.Range("M77").Value = -3779.5
$ws.Range("N77").Value = -30464.845
$ws.Range("H102").Value = 10041193
$ws.Range("I102").Value = 1263.625
$ws.Range("J102").Value = 27889956
$ws.Range("K102").Value = 1263.625
$ws.Range("L102").Value = 27889956
$ws.Range("M102").Value = 358.375
$ws.Range("N102").Value = -27893200
$ws.Range("H132").Value = 4214.136
$ws.Range("I132").Value = 3271.5193
$ws.Range("J132").Value = 7715.2856
$ws.Range("K132").Value = 9814.5579
$ws.Range("L132").Value = 23145.8568
$ws.Range("M132").Value = -7284.5579
$ws.Range("N132").Value = -28205.8568
$ws.Range("H136").Value = 5411.857
$ws.Range("I136").Value = 4590.974
$ws.Range("K136").Value = 13772.922
$ws.Range("M136").Value = -11222.922

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1988.6
$ws.Range("I22").Value = 1988.6
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1988.6
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1815.6
$ws.Range("N22").ClearContents()
$ws.Range("H32").Value = 60026.5
$ws.Range("J32").Value = 60026.5
$ws.Range("L32").Value = 60026.5
$ws.Range("N32").Value = -60794.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2862012.5
$ws.Range("I13").Value = 20000000
$ws.Range("J13").Value = 5681.3335
$ws.Range("K13").Value = 20000000
$ws.Range("L13").Value = 5681.3335
$ws.Range("M13").Value = -19999861
$ws.Range("N13").Value = -5959.3335
$ws.Range("H14").Value = 1344.4667
$ws.Range("I14").Value = 2074.75
$ws.Range("J14").Value = 1078.909
$ws.Range("K14").Value = 2074.75
$ws.Range("L14").Value = 1078.909
$ws.Range("M14").Value = -1904.75
$ws.Range("N14").Value = -1418.909
$ws.Range("H21").Value = 9993
$ws.Range("J21").Value = 9993
$ws.Range("L21").Value = 9993
$ws.Range("N21").Value = -10463
$ws.Range("H26").Value = 4000
$ws.Range("J26").Value = 4000
$ws.Range("L26").Value = 4000
$ws.Range("N26").Value = -4574
$ws.Range("H31").Value = 22226268
$ws.Range("J31").Value = 5221.8335
$ws.Range("L31").Value = 5221.8335
$ws.Range("N31").Value = -5811.8335
$ws.Range("H34").Value = 22226268
$ws.Range("J34").Value = 5221.8335
$ws.Range("L34").Value = 5221.8335
$ws.Range("N34").Value = -5625.8335
$ws.Range("H99").Value = 9861.75
$ws.Range("I99").Value = 8749
$ws.Range("K99").Value = 8749
$ws.Range("M99").Value = -7251
$ws.Range("H100").Value = 55375
$ws.Range("J100").Value = 55375
$ws.Range("L100").Value = 55375
$ws.Range("N100").Value = -57539
$ws.Range("H105").Value = 1295.7142
$ws.Range("I105").Value = 1157.2222
$ws.Range("K105").Value = 1157.2222
$ws.Range("M105").Value = 589.7778000000001
$ws.Range("H126").Value = 9861.75
$ws.Range("I126").Value = 8749
$ws.Range("K126").Value = 26247
$ws.Range("M126").Value = -23777
$ws.Range("H132").Value = 2671.8948
$ws.Range("I132").Value = 1827.8667
$ws.Range("K132").Value = 5483.6001
$ws.Range("M132").Value = -2953.6001
$ws.Range("H134").Value = 7446.8047
$ws.Range("I134").Value = 7080.3438
$ws.Range("J134").Value = 8749.777
$ws.Range("K134").Value = 21241.0314
$ws.Range("L134").Value = 26249.331
$ws.Range("M134").Value = -18706.0314
$ws.Range("N134").Value = -31319.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3354.0476
$ws.Range("I2").Value = 46.6
$ws.Range("K2").Value = 279.6
$ws.Range("M2").Value = -166.6
$ws.Range("H4").Value = 3339869
$ws.Range("I4").Value = 2584754.2
$ws.Range("K4").Value = 7754262.600000001
$ws.Range("M4").Value = -7754150.600000001
$ws.Range("H13").Value = 234.42857
$ws.Range("J13").Value = 270
$ws.Range("L13").Value = 810
$ws.Range("N13").Value = -1146
$ws.Range("H122").Value = 166678850
$ws.Range("I122").Value = 333356740
$ws.Range("J122").Value = 961.6667
$ws.Range("K122").Value = 3000210660
$ws.Range("L122").Value = 8655.0003
$ws.Range("M122").Value = -3000208210
$ws.Range("N122").Value = -13555.0003
$ws.Range("H131").Value = 5378.0557
$ws.Range("I131").Value = 1615
$ws.Range("J131").Value = 5599.4116
$ws.Range("K131").Value = 4845
$ws.Range("L131").Value = 16798.2348
$ws.Range("M131").Value = 195
$ws.Range("N131").Value = -26878.2348
$ws.Range("H132").Value = 4907
$ws.Range("I132").Value = 6098.3335
$ws.Range("J132").Value = 1333
$ws.Range("K132").Value = 54885.0015
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -52355.0015
$ws.Range("N132").Value = -17057
$ws.Range("H134").Value = 38463668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 12911.467
$ws.Range("I113").Value = 15038.333
$ws.Range("J113").Value = 4404
$ws.Range("K113").Value = 15038.333
$ws.Range("L113").Value = 4404
$ws.Range("M113").Value = -12868.333
$ws.Range("N113").Value = -8744
$ws.Range("H132").Value = 3895.68
$ws.Range("I132").Value = 2501.5
$ws.Range("K132").Value = 7504.5
$ws.Range("M132").Value = -4974.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 479.26315
$ws.Range("I16").Value = 479.26315
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 479.26315
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -309.26315
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 2770.1177
$ws.Range("I22").Value = 1712.28
$ws.Range("J22").Value = 3787.2693
$ws.Range("K22").Value = 1712.28
$ws.Range("L22").Value = 3787.2693
$ws.Range("M22").Value = -1417.28
$ws.Range("N22").Value = -4377.2693
$ws.Range("H27").Value = 2770.1177
$ws.Range("I27").Value = 1712.28
$ws.Range("J27").Value = 3787.2693
$ws.Range("K27").Value = 1712.28
$ws.Range("L27").Value = 3787.2693
$ws.Range("M27").Value = -1605.28
$ws.Range("N27").Value = -4001.2693
$ws.Range("H40").Value = 26323780
$ws.Range("I40").Value = 31256770
$ws.Range("K40").Value = 31256770
$ws.Range("M40").Value = -31256634
$ws.Range("H46").Value = 4383.204
$ws.Range("I46").Value = 1700
$ws.Range("J46").Value = 5071.205
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 5071.205
$ws.Range("M46").Value = -1512
$ws.Range("N46").Value = -5447.205
$ws.Range("H132").Value = 2240.0488
$ws.Range("J132").Value = 9636
$ws.Range("L132").Value = 28908
$ws.Range("N132").Value = -33968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17474.285
$ws.Range("I62").Value = 17553
$ws.Range("J62").Value = 17002
$ws.Range("K62").Value = 17553
$ws.Range("L62").Value = 17002
$ws.Range("M62").Value = -16929
$ws.Range("N62").Value = -18250
$ws.Range("H65").Value = 17474.285
$ws.Range("I65").Value = 17553
$ws.Range("J65").Value = 17002
$ws.Range("K65").Value = 87765
$ws.Range("L65").Value = 85010
$ws.Range("M65").Value = -84645
$ws.Range("N65").Value = -91250
$ws.Range("H96").Value = 18724.75
$ws.Range("I96").Value = 8450
$ws.Range("K96").Value = 8450
$ws.Range("M96").Value = -7077
$ws.Range("H113").Value = 940.20514
$ws.Range("I113").Value = 553
$ws.Range("K113").Value = 1659
$ws.Range("M113").Value = 511
$ws.Range("H126").Value = 5113.8945
$ws.Range("I126").Value = 5321.615
$ws.Range("J126").Value = 4663.8335
$ws.Range("K126").Value = 15964.845
$ws.Range("L126").Value = 13991.5005
$ws.Range("M126").Value = -13494.845
$ws.Range("N126").Value = -18931.5005
$ws.Range("H132").Value = 3155.1929
$ws.Range("I132").Value = 2167.4146
$ws.Range("K132").Value = 6502.2438
$ws.Range("M132").Value = -3972.2438
